$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$objeto = "Contratação de pessoa jurídica para fornecimento de material elétrico para ILUMINAÇÃO PÚBLICA em atendimento às necessidades da Secretaria de Infraestrutura e Urbanismo do Município de Nilo Peçanha - BA., na forma estabelecida no Termo de Referência e de acordo com a proposta do contratado que para todos os efeitos integra este contrato como se transcrita fosse, apresentada na forma de anexo único ao presente."

foreach ($row in 76, 77) {
    # Columns A:J and U stay blank, but the cells still need to exist in sheetData.
    $ws.Range("A$row`:J$row").Font.Bold = $false
    $ws.Range("U$row").Font.Bold = $false

    $ws.Range("K$row").Value = "154/2023"
    $ws.Range("L$row").Value = "068/2023"
    $ws.Range("M$row").Value = "19/12/2023"
    $ws.Range("N$row").Value = "31/12/2023"
    $ws.Range("O$row").Value = "LABORCOM COMÉRCIO DE MATERIAIS DE CONSTRUÇÃO LTDA."
    $ws.Range("P$row").Value = "34.101.659/0001-56"
    $ws.Range("Q$row").Value = "DISPENSA DE LICITAÇÃO"
    $ws.Range("R$row").Value = $objeto
    $ws.Range("S$row").Value = "MUNICÍPIO DE NILO PEÇANHA"
    $ws.Range("T$row").Value = "54.720,00"
}
